$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder shared data: ECs now referenced before sCs, Cdh1, Egfr stay, FAPs stays.
# Rebuild rows 2-7 with sending/target cluster combinations per updated NATMI output.

$rows = @(
    @{ A="ECs"; B="Cdh1"; C="Egfr"; D="ECs";  E=1; F=0.3333333333333333; G=0.098866; H=0.296598; I=0.3026185969870575; J=0.3026185969870575; K=3; L=1; M=1.307106666666667;  N=3.92132;         O=0.01256263154946851; P=0.01256263154946851; Q=0.1292284077066667; R=1.16305566936;     S=0.003801685933965505; T=0.003801685933965504 },
    @{ A="ECs"; B="Cdh1"; C="Egfr"; D="FAPs"; E=1; F=0.3333333333333333; G=0.098866; H=0.296598; I=0.3026185969870575; J=0.3026185969870575; K=3; L=1; M=80.22623699999998;  N=240.678711;      O=0.77105616682495;    P=0.77105616682495;    Q=7.931647147241998;  R=71.38482432517799; S=0.2333359354027849;   T=0.2333359354027849 },
    @{ A="ECs"; B="Cdh1"; C="Egfr"; D="sCs";  E=1; F=0.3333333333333333; G=0.098866; H=0.296598; I=0.3026185969870575; J=0.3026185969870575; K=3; L=1; M=22.51385866666667;  N=67.54157600000001; O=0.2163812016255815;  P=0.2163812016255815;  Q=2.225855150938667;  R=20.032696358448;   S=0.06548097565030707;  T=0.06548097565030706 },
    @{ A="sCs"; B="Cdh1"; C="Egfr"; D="ECs";  E=3; F=1;                  G=0.2278356666666667; H=0.6835070000000001; I=0.6973814030129426; J=0.6973814030129426; K=3; L=1; M=1.307106666666667;  N=3.92132;         O=0.01256263154946851; P=0.01256263154946851; Q=0.2978055188044444; R=2.680249669240001; S=0.008760945615503007; T=0.008760945615503007 },
    @{ A="sCs"; B="Cdh1"; C="Egfr"; D="FAPs"; E=3; F=1;                  G=0.2278356666666667; H=0.6835070000000001; I=0.6973814030129426; J=0.6973814030129426; K=3; L=1; M=80.22623699999998;  N=240.678711;      O=0.77105616682495;    P=0.77105616682495;    Q=18.278398191053;    R=164.505583719477;  S=0.5377202314221652;   T=0.5377202314221652 },
    @{ A="sCs"; B="Cdh1"; C="Egfr"; D="sCs";  E=3; F=1;                  G=0.2278356666666667; H=0.6835070000000001; I=0.6973814030129426; J=0.6973814030129426; K=3; L=1; M=22.51385866666667;  N=67.54157600000001; O=0.2163812016255815;  P=0.2163812016255815;  Q=5.129459998559112;  R=46.16513998703201; S=0.1509002259752744;   T=0.1509002259752744 }
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $rows[$i]
    foreach ($col in $columns) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
